$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.487.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.561.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.98%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.39%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'211.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.32%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.87%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.38%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'45.97"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +5.49%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'24.06"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.85%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -2.04%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -1.81%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -0.53%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D14").Value = "'1.570.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.27%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -2.27%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'28.488.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.11%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -3.04%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'61.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -3.36%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'226.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.02%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -2.26%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -3.18%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.42%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -7.12%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'9.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.08%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +6.52%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'149.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.17%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -2.69%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -3.14%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -3.27%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.38%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "'Hedera"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'0.0464"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.20%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "'PancakeSwap"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'1.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.75%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.43%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.52%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.394.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.01%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.96%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -4.53%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +1.46%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.84%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.0166"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.31%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -2.00%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +0.40%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.786"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.63%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.96%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'5.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.83%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.978"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.52%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'62.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.13%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.698.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.80%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'86.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.00%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0₆0103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.16%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0518"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.28%  "
$ws.Range("E51").Style = "Normal"
